$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Status changed from "In Progress" to "Done" ---
$ws.Range("J8").Value = "Done"
$ws.Range("J8").Interior.Color = 65280      # 00FF00 - "Done" status colour

# --- Row 9: Status changed from "Done" to "Not Done" ---
$ws.Range("J9").Value = "Not Done"
$ws.Range("J9").Interior.Color = 255        # FF0000 - "Not Done" status colour

# --- Row 15: Priority/Risk/Points filled in (previously "TBD") ---
$ws.Range("G15").Value = "Could"
$ws.Range("H15").Value = "Low/High"
$ws.Range("I15").Value = 1

# --- Row 16: Priority/Risk/Points filled in (previously "TBD") ---
$ws.Range("G16").Value = "Could"
$ws.Range("H16").Value = "High/Low"
$ws.Range("I16").Value = 3

# --- Row 17: Priority/Risk/Points filled in (previously "TBD") ---
$ws.Range("G17").Value = "Could"
$ws.Range("H17").Value = "High/Low"
$ws.Range("I17").Value = "?"

# --- Row 18: brand-new backlog item (User Story 24) ---
$ws.Range("C18").Value = 24
$ws.Range("D18").Value = "Site Visitor"
$ws.Range("E18").Value = "Access a map of events that are related to a goal that I am interested in"
$ws.Range("F18").Value = "See which events are in my area"
$ws.Range("G18").Value = "Must"
$ws.Range("H18").Value = "High/High"
$ws.Range("I18").Value = 13
$ws.Range("J18").Value = "Not Done"

# Match the established colour-coding used throughout the sheet for each column
$ws.Range("C18").Interior.Color = 2441676   # CC4125 - matches other "User Story ID" cells
$ws.Range("D18").Interior.Color = 8242323   # 93C47D - matches other "User" cells
$ws.Range("E18").Interior.Color = 6740479   # FFD966 - matches other "I want to..." cells
$ws.Range("F18").Interior.Color = 7949734   # A64D79 - matches other "So that I can..." cells
$ws.Range("G18").Interior.Color = 2441676   # CC4125 - matches other "Priority" cells
$ws.Range("H18").Interior.Color = 5482548   # 34A853 - matches other "Risk/Value" cells
$ws.Range("I18").Interior.Color = 310523    # FBBC04 - matches other "Points" cells
$ws.Range("J18").Interior.Color = 255       # FF0000 - matches other "Not Done" status cells
